{"js": "// The document contains one table of three-digit-number \u00f7 one-digit-number\n// division problems (\"dividend\u00f7divisor=quotient, remainder\"). Several rows\n// are blank spacer rows. This script replaces the text of 25 specific\n// table cells (identified by row/column index, since some of the original\n// values repeat verbatim elsewhere in the table) with new problem/answer\n// strings, while preserving each cell's existing run/paragraph formatting.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"373\u00f75=74, 3\", newText: \"625\u00f79=69, 4\" },\n  { row: 0, col: 1, oldText: \"995\u00f79=110, 5\", newText: \"484\u00f76=80, 4\" },\n  { row: 0, col: 2, oldText: \"276\u00f74=69, 0\", newText: \"303\u00f75=60, 3\" },\n  { row: 0, col: 3, oldText: \"176\u00f73=58, 2\", newText: \"655\u00f78=81, 7\" },\n  { row: 0, col: 4, oldText: \"408\u00f72=204, 0\", newText: \"265\u00f73=88, 1\" },\n  { row: 4, col: 0, oldText: \"116\u00f72=58, 0\", newText: \"781\u00f75=156, 1\" },\n  { row: 4, col: 1, oldText: \"869\u00f77=124, 1\", newText: \"176\u00f76=29, 2\" },\n  { row: 4, col: 2, oldText: \"441\u00f72=220, 1\", newText: \"375\u00f77=53, 4\" },\n  { row: 4, col: 3, oldText: \"631\u00f75=126, 1\", newText: \"714\u00f72=357, 0\" },\n  { row: 4, col: 4, oldText: \"619\u00f74=154, 3\", newText: \"612\u00f77=87, 3\" },\n  { row: 8, col: 0, oldText: \"826\u00f77=118, 0\", newText: \"702\u00f73=234, 0\" },\n  { row: 8, col: 1, oldText: \"164\u00f78=20, 4\", newText: \"904\u00f77=129, 1\" },\n  { row: 8, col: 2, oldText: \"453\u00f73=151, 0\", newText: \"939\u00f72=469, 1\" },\n  { row: 8, col: 3, oldText: \"985\u00f76=164, 1\", newText: \"354\u00f74=88, 2\" },\n  { row: 8, col: 4, oldText: \"619\u00f74=154, 3\", newText: \"271\u00f76=45, 1\" },\n  { row: 12, col: 0, oldText: \"315\u00f77=45, 0\", newText: \"632\u00f78=79, 0\" },\n  { row: 12, col: 1, oldText: \"163\u00f74=40, 3\", newText: \"846\u00f77=120, 6\" },\n  { row: 12, col: 2, oldText: \"768\u00f78=96, 0\", newText: \"328\u00f73=109, 1\" },\n  { row: 12, col: 3, oldText: \"748\u00f77=106, 6\", newText: \"520\u00f78=65, 0\" },\n  { row: 12, col: 4, oldText: \"124\u00f77=17, 5\", newText: \"880\u00f77=125, 5\" },\n  { row: 16, col: 0, oldText: \"146\u00f78=18, 2\", newText: \"915\u00f79=101, 6\" },\n  { row: 16, col: 1, oldText: \"720\u00f77=102, 6\", newText: \"380\u00f75=76, 0\" },\n  { row: 16, col: 2, oldText: \"953\u00f74=238, 1\", newText: \"394\u00f72=197, 0\" },\n  { row: 16, col: 3, oldText: \"666\u00f79=74, 0\", newText: \"426\u00f77=60, 6\" },\n  { row: 16, col: 4, oldText: \"567\u00f72=283, 1\", newText: \"142\u00f78=17, 6\" }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const range = cell.body.getRange();\n  range.load(\"text\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  // Word table cell ranges report their text with a trailing end-of-cell\n  // marker (\\r or \\t); strip it before comparing.\n  const actual = range.text.replace(/[\\t\\r]+$/, \"\");\n  if (actual !== r.oldText) {\n    throw new Error(\n      \"Unexpected cell text at row \" + r.row + \", col \" + r.col +\n      \": got \" + JSON.stringify(range.text) + \", expected \" + JSON.stringify(r.oldText)\n    );\n  }\n\n  range.insertText(r.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document's single table holds three-digit-number / one-digit-number\n# division problems (\"dividend\u00f7divisor=quotient, remainder\"); 4 out of\n# every 5 rows are blank spacer rows. This script overwrites the text of\n# 25 specific table cells (addressed by 1-based row/column, since some\n# original values repeat verbatim elsewhere in the table) with new\n# problem/answer strings while leaving each cell's existing paragraph and\n# run formatting untouched (Range.Text only replaces the text run(s)).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$items = @(\n    @{ Row = 1; Col = 1; Old = \"373\u00f75=74, 3\"; New = \"625\u00f79=69, 4\" },\n    @{ Row = 1; Col = 2; Old = \"995\u00f79=110, 5\"; New = \"484\u00f76=80, 4\" },\n    @{ Row = 1; Col = 3; Old = \"276\u00f74=69, 0\"; New = \"303\u00f75=60, 3\" },\n    @{ Row = 1; Col = 4; Old = \"176\u00f73=58, 2\"; New = \"655\u00f78=81, 7\" },\n    @{ Row = 1; Col = 5; Old = \"408\u00f72=204, 0\"; New = \"265\u00f73=88, 1\" },\n    @{ Row = 5; Col = 1; Old = \"116\u00f72=58, 0\"; New = \"781\u00f75=156, 1\" },\n    @{ Row = 5; Col = 2; Old = \"869\u00f77=124, 1\"; New = \"176\u00f76=29, 2\" },\n    @{ Row = 5; Col = 3; Old = \"441\u00f72=220, 1\"; New = \"375\u00f77=53, 4\" },\n    @{ Row = 5; Col = 4; Old = \"631\u00f75=126, 1\"; New = \"714\u00f72=357, 0\" },\n    @{ Row = 5; Col = 5; Old = \"619\u00f74=154, 3\"; New = \"612\u00f77=87, 3\" },\n    @{ Row = 9; Col = 1; Old = \"826\u00f77=118, 0\"; New = \"702\u00f73=234, 0\" },\n    @{ Row = 9; Col = 2; Old = \"164\u00f78=20, 4\"; New = \"904\u00f77=129, 1\" },\n    @{ Row = 9; Col = 3; Old = \"453\u00f73=151, 0\"; New = \"939\u00f72=469, 1\" },\n    @{ Row = 9; Col = 4; Old = \"985\u00f76=164, 1\"; New = \"354\u00f74=88, 2\" },\n    @{ Row = 9; Col = 5; Old = \"619\u00f74=154, 3\"; New = \"271\u00f76=45, 1\" },\n    @{ Row = 13; Col = 1; Old = \"315\u00f77=45, 0\"; New = \"632\u00f78=79, 0\" },\n    @{ Row = 13; Col = 2; Old = \"163\u00f74=40, 3\"; New = \"846\u00f77=120, 6\" },\n    @{ Row = 13; Col = 3; Old = \"768\u00f78=96, 0\"; New = \"328\u00f73=109, 1\" },\n    @{ Row = 13; Col = 4; Old = \"748\u00f77=106, 6\"; New = \"520\u00f78=65, 0\" },\n    @{ Row = 13; Col = 5; Old = \"124\u00f77=17, 5\"; New = \"880\u00f77=125, 5\" },\n    @{ Row = 17; Col = 1; Old = \"146\u00f78=18, 2\"; New = \"915\u00f79=101, 6\" },\n    @{ Row = 17; Col = 2; Old = \"720\u00f77=102, 6\"; New = \"380\u00f75=76, 0\" },\n    @{ Row = 17; Col = 3; Old = \"953\u00f74=238, 1\"; New = \"394\u00f72=197, 0\" },\n    @{ Row = 17; Col = 4; Old = \"666\u00f79=74, 0\"; New = \"426\u00f77=60, 6\" },\n    @{ Row = 17; Col = 5; Old = \"567\u00f72=283, 1\"; New = \"142\u00f78=17, 6\" }\n)\n\nforeach ($item in $items) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    # Cell Range.Text includes the trailing end-of-cell marker (CR + BEL);\n    # strip it before comparing against the expected original text.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $item.Old) {\n        throw \"Unexpected text at row $($item.Row), col $($item.Col): got [$current], expected [$($item.Old)]\"\n    }\n    $cell.Range.Text = $item.New\n}\n"}
